$d = $word.ActiveDocument

$d.Content.Find.Execute("113÷9=12, 5", $true, $false, $false, $false, $false, $true, 1, $false, "970÷7=138, 4", 2) | Out-Null
$d.Content.Find.Execute("344÷6=57, 2", $true, $false, $false, $false, $false, $true, 1, $false, "535÷6=89, 1", 2) | Out-Null
$d.Content.Find.Execute("724÷7=103, 3", $true, $false, $false, $false, $false, $true, 1, $false, "432÷8=54, 0", 2) | Out-Null
$d.Content.Find.Execute("960÷3=320, 0", $true, $false, $false, $false, $false, $true, 1, $false, "256÷8=32, 0", 2) | Out-Null
$d.Content.Find.Execute("204÷8=25, 4", $true, $false, $false, $false, $false, $true, 1, $false, "984÷5=196, 4", 2) | Out-Null
$d.Content.Find.Execute("823÷2=411, 1", $true, $false, $false, $false, $false, $true, 1, $false, "675÷8=84, 3", 2) | Out-Null
$d.Content.Find.Execute("216÷2=108, 0", $true, $false, $false, $false, $false, $true, 1, $false, "328÷2=164, 0", 2) | Out-Null
$d.Content.Find.Execute("342÷8=42, 6", $true, $false, $false, $false, $false, $true, 1, $false, "346÷6=57, 4", 2) | Out-Null
$d.Content.Find.Execute("809÷5=161, 4", $true, $false, $false, $false, $false, $true, 1, $false, "196÷8=24, 4", 2) | Out-Null
$d.Content.Find.Execute("639÷7=91, 2", $true, $false, $false, $false, $false, $true, 1, $false, "917÷9=101, 8", 2) | Out-Null
$d.Content.Find.Execute("556÷4=139, 0", $true, $false, $false, $false, $false, $true, 1, $false, "394÷9=43, 7", 2) | Out-Null
$d.Content.Find.Execute("903÷2=451, 1", $true, $false, $false, $false, $false, $true, 1, $false, "524÷8=65, 4", 2) | Out-Null
$d.Content.Find.Execute("394÷8=49, 2", $true, $false, $false, $false, $false, $true, 1, $false, "266÷3=88, 2", 2) | Out-Null
$d.Content.Find.Execute("410÷3=136, 2", $true, $false, $false, $false, $false, $true, 1, $false, "628÷4=157, 0", 2) | Out-Null
$d.Content.Find.Execute("506÷8=63, 2", $true, $false, $false, $false, $false, $true, 1, $false, "783÷4=195, 3", 2) | Out-Null
$d.Content.Find.Execute("163÷4=40, 3", $true, $false, $false, $false, $false, $true, 1, $false, "125÷9=13, 8", 2) | Out-Null
$d.Content.Find.Execute("452÷4=113, 0", $true, $false, $false, $false, $false, $true, 1, $false, "739÷2=369, 1", 2) | Out-Null
$d.Content.Find.Execute("745÷5=149, 0", $true, $false, $false, $false, $false, $true, 1, $false, "670÷6=111, 4", 2) | Out-Null
$d.Content.Find.Execute("531÷2=265, 1", $true, $false, $false, $false, $false, $true, 1, $false, "917÷4=229, 1", 2) | Out-Null
$d.Content.Find.Execute("484÷3=161, 1", $true, $false, $false, $false, $false, $true, 1, $false, "708÷8=88, 4", 2) | Out-Null
$d.Content.Find.Execute("563÷2=281, 1", $true, $false, $false, $false, $false, $true, 1, $false, "969÷3=323, 0", 2) | Out-Null
$d.Content.Find.Execute("267÷8=33, 3", $true, $false, $false, $false, $false, $true, 1, $false, "497÷6=82, 5", 2) | Out-Null
$d.Content.Find.Execute("214÷6=35, 4", $true, $false, $false, $false, $false, $true, 1, $false, "701÷2=350, 1", 2) | Out-Null
$d.Content.Find.Execute("168÷9=18, 6", $true, $false, $false, $false, $false, $true, 1, $false, "212÷2=106, 0", 2) | Out-Null
$d.Content.Find.Execute("748÷2=374, 0", $true, $false, $false, $false, $false, $true, 1, $false, "857÷6=142, 5", 2) | Out-Null
